$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure cells keep their existing text type (no implicit numeric/date conversion)
$priceCells = @("D2","D3","D5","D6","D7","D13","D16","D17","D18","D19","D20","D21","D22","D24","D25","D29","D30","D32","D33","D36","D39","D40","D41","D42","D43","D49","D50","D51")
$volCells = @("E2","E3","E4","E5","E6","E7","E8","E9","E10","E11","E12","E13","E14","E15","E16","E17","E18","E19","E20","E21","E22","E23","E24","E25","E26","E27","E28","E29","E30","E31","E32","E33","E34","E35","E36","E37","E38","E39","E40","E41","E42","E43","E44","E45","E46","E47","E48","E49","E50","E51")
foreach ($cellRef in $priceCells) { $ws.Range($cellRef).NumberFormat = "@" }
foreach ($cellRef in $volCells) { $ws.Range($cellRef).NumberFormat = "@" }

# Row 2 - Bitcoin
$ws.Range("D2").Value = "64.459.31"
$ws.Range("E2").Value = "  -0.13%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.495.99"
$ws.Range("E3").Value = "  -0.02%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.11%  "

# Row 5 - BNB
$ws.Range("D5").Value = "586.85"
$ws.Range("E5").Value = "  +0.30%  "

# Row 6 - Solana
$ws.Range("D6").Value = "135.06"
$ws.Range("E6").Value = "  +2.67%  "

# Row 7 - LidoStakedEther
$ws.Range("D7").Value = "3.497.95"
$ws.Range("E7").Value = "  +0.16%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  -0.03%  "

# Row 9 - XRP
$ws.Range("E9").Value = "  -0.78%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  +0.11%  "

# Row 11 - Toncoin
$ws.Range("E11").Value = "  -0.76%  "

# Row 12 - Cardano
$ws.Range("E12").Value = "  -2.84%  "

# Row 13 - WrappedliquidstakedEther2.0
$ws.Range("D13").Value = "4.091.75"
$ws.Range("E13").Value = "  -0.28%  "

# Row 14 - ShibaInu
$ws.Range("E14").Value = "  +0.83%  "

# Row 15 - TRON
$ws.Range("E15").Value = "  +1.33%  "

# Row 16 - WrappedEther
$ws.Range("D16").Value = "3.494.79"
$ws.Range("E16").Value = "  -0.31%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "64.442.48"
$ws.Range("E17").Value = "  -0.23%  "

# Row 18 - Avalanche
$ws.Range("D18").Value = "25.32"
$ws.Range("E18").Value = "  -8.43%  "

# Row 19 - Uniswap
$ws.Range("D19").Value = "10.00"
$ws.Range("E19").Value = "  +0.94%  "

# Row 20 - Polkadot
$ws.Range("D20").Value = "5.64"
$ws.Range("E20").Value = "  -0.15%  "

# Row 21 - Chainlink
$ws.Range("D21").Value = "13.77"
$ws.Range("E21").Value = "  -4.33%  "

# Row 22 - BitcoinCash
$ws.Range("D22").Value = "385.94"
$ws.Range("E22").Value = "  -1.89%  "

# Row 23 - Polygon
$ws.Range("E23").Value = "  -1.45%  "

# Row 24 - WrappedeETH
$ws.Range("D24").Value = "3.634.75"
$ws.Range("E24").Value = "  -0.27%  "

# Row 25 - Litecoin
$ws.Range("D25").Value = "74.07"
$ws.Range("E25").Value = "  +0.90%  "

# Row 26 - Dai
$ws.Range("E26").Value = "  +0.07%  "

# Row 27 - PEPE
$ws.Range("E27").Value = "  +4.05%  "

# Row 28 - RenderToken
$ws.Range("E28").Value = "  +0.39%  "

# Row 29 - was Binance-PegBSC-USD, now Fetch.AI
$ws.Range("B29").Value = "Fetch.AI"
$ws.Range("C29").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D29").Value = "1.54"
$ws.Range("E29").Value = "  -1.88%  "

# Row 30 - was Fetch.AI, now Binance-PegBSC-USD
$ws.Range("B30").Value = "Binance-PegBSC-USD"
$ws.Range("C30").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  +0.05%  "

# Row 31 - PancakeSwap
$ws.Range("E31").Value = "  -1.07%  "

# Row 32 - InternetComputer(DFINITY)
$ws.Range("D32").Value = "8.23"
$ws.Range("E32").Value = "  +0.67%  "

# Row 33 - RenzoRestakedETH
$ws.Range("D33").Value = "3.519.58"
$ws.Range("E33").Value = "  +0.44%  "

# Row 34 - USDe
$ws.Range("E34").Value = "  -0.01%  "

# Row 35 - Kaspa
$ws.Range("E35").Value = "  +0.29%  "

# Row 36 - EthereumClassic
$ws.Range("D36").Value = "23.46"
$ws.Range("E36").Value = "  -2.47%  "

# Row 37 - NEARProtocol
$ws.Range("E37").Value = "  +1.61%  "

# Row 38 - ImmutableX
$ws.Range("E38").Value = "  -2.80%  "

# Row 39 - Aptos
$ws.Range("D39").Value = "6.84"
$ws.Range("E39").Value = "  -1.76%  "

# Row 40 - Monero
$ws.Range("D40").Value = "162.44"
$ws.Range("E40").Value = "  -4.92%  "

# Row 41 - Hedera
$ws.Range("D41").Value = "0.0780"
$ws.Range("E41").Value = "  -2.87%  "

# Row 42 - Mantle
$ws.Range("D42").Value = "0.804"
$ws.Range("E42").Value = "  -1.06%  "

# Row 43 - EnergySwap
$ws.Range("D43").Value = "25.90"
$ws.Range("E43").Value = "  -1.04%  "

# Row 44 - FirstDigitalUSD
$ws.Range("E44").Value = "  -0.14%  "

# Row 45 - OKB
$ws.Range("E45").Value = "  -0.07%  "

# Row 46 - ONDO
$ws.Range("E46").Value = "  +0.75%  "

# Row 47 - Filecoin
$ws.Range("E47").Value = "  +1.27%  "

# Row 48 - Stacks
$ws.Range("E48").Value = "  +0.55%  "

# Row 49 - Maker
$ws.Range("D49").Value = "2.478.31"
$ws.Range("E49").Value = "  +1.64%  "

# Row 50 - Cosmos
$ws.Range("D50").Value = "6.76"
$ws.Range("E50").Value = "  -1.53%  "

# Row 51 - SuiNetwork
$ws.Range("D51").Value = "0.906"
$ws.Range("E51").Value = "  +1.85%  "
